# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Terminal Hortofrutícola Agro Chillán - Mango"
# at row 173, pushing the existing rows 173:211 down to 174:212, then fill in
# the new row with the latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 173 (shifts 173:211 -> 174:212).
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with this week's record.
$ws.Range("A173").Value = 7
$ws.Range("B173").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C173").Value = "Ñuble"
$ws.Range("D173").Value = 45244
$ws.Range("E173").Value = 16
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100108
$ws.Range("H173").Value = "Tropicales y subtropicales"
$ws.Range("I173").Value = 100108002
$ws.Range("J173").Value = "Mango"
$ws.Range("K173").Value = "Sin especificar"
$ws.Range("L173").Value = "Primera"
$ws.Range("M173").Value = 60
$ws.Range("N173").Value = 12000
$ws.Range("O173").Value = 12000
$ws.Range("P173").Value = 12000
$ws.Range("Q173").Value = '$/bandeja 4 kilos'
$ws.Range("R173").Value = "Brasil"
$ws.Range("S173").Value = 3000
$ws.Range("T173").Value = 4
